$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (fix bug in fedrollover): rows 390-393
$newData = @(
    @(20082300, 0),
    @(20082400, 0),
    @(20082500, 6218625100),
    @(20082600, 0)
)

$startRow = 390
for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newData[$i][0]
    $ws.Cells.Item($r, 2).Value = $newData[$i][1]
}
